# Adds a new "Pre Conditions" column (C) to the Test Cases sheet, between
# "Test Scenario" and "Test Steps", and fills it with data for every test
# case row. Everything that used to live in columns C:G shifts right to D:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# --- Insert the new column -------------------------------------------------
# Inserting a whole column shifts the existing C:G data/styles to D:H and
# keeps A:B untouched - exactly the shape described by the diff.
$ws.Columns("C:C").Insert()

# --- Column width: new column C should be as wide as column B -------------
$bWidth = $ws.Columns("B:B").ColumnWidth
$ws.Columns("C:C").ColumnWidth = $bWidth

# --- Header row -------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial($xlPasteFormats)
$ws.Range("C1").Value = "Pre Conditions"

# --- Data rows ---------------------------------------------------------------
# Style "4" (center/center, no wrap) is used for most Pre Conditions cells;
# style "6" (center/center, wrap) is used for the couple of rows whose text
# is longer. Grab untouched donor cells that already carry those styles so
# the copy reuses the existing style indexes instead of synthesizing new
# ones.
$ws.Range("G2").Copy()
$style4Donor = "G2"
$ws.Range("F3").Copy()
$style6Donor = "F3"

function Set-PreCondition($row, $text, $style) {
    if ($style -eq 6) {
        $ws.Range($style6Donor).Copy()
    } else {
        $ws.Range($style4Donor).Copy()
    }
    $ws.Range("C$row").PasteSpecial($xlPasteFormats)
    $ws.Range("C$row").Value = $text
}

$noControlChange = "Ustawienia domyślne sterowania `nnie zostały zmienione."
$notLastLevel = "Poziom rozpoczęcia gry nie jest `nostanim poziomem gry"
$dash = "-"

Set-PreCondition 2  $noControlChange 6
Set-PreCondition 3  $dash 4
Set-PreCondition 4  $dash 4
Set-PreCondition 5  $dash 4
Set-PreCondition 6  $dash 4
Set-PreCondition 7  $dash 4
Set-PreCondition 8  $dash 4
Set-PreCondition 9  $dash 4
Set-PreCondition 10 $dash 4
Set-PreCondition 11 $dash 4
Set-PreCondition 12 $dash 4
Set-PreCondition 13 $dash 4
Set-PreCondition 14 $dash 4
Set-PreCondition 15 $dash 4
Set-PreCondition 16 $notLastLevel 6
Set-PreCondition 17 $dash 4
Set-PreCondition 18 $dash 4
Set-PreCondition 19 $notLastLevel 6
Set-PreCondition 20 $dash 4
Set-PreCondition 21 $dash 4
Set-PreCondition 22 $dash 4
Set-PreCondition 23 $dash 4
Set-PreCondition 24 $dash 4

# --- View: zoom + selected cell ---------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("M5").Select()
